# Modification de la taille des images
# (rows describing image-size related SEO issues are removed,
#  and the category of the remaining "missing title tags" row is
#  changed to a new "mauvaise pratique html" category)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows that no longer apply:
#  - row 22: "class et id dans la même div"
#  - row 24: "nom des class et id"
#  - row 30: "alt avant src"
$ws.Range("A22:E22").ClearContents()
$ws.Range("A24:E24").ClearContents()
$ws.Range("A30:E30").ClearContents()

# Row 25 ("balise sans texte") is recategorised from "seo" to the
# new category "mauvaise pratique html"
$ws.Range("A25").Value = "mauvaise pratique html"

# Update the current selection/cursor position in the sheet view
$ws.Range("A34:A35").Select()
